$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT into a cell, preserving the cell's
# original (default) style. Values in the Price/Volume columns look like
# numbers (e.g. "30.646.06", "0.9997") or percentages, and Excel's normal
# Range.Value setter would silently parse/convert them into real numbers,
# which does not match the source (text) cells. Flipping the number format
# to Text before the write keeps the literal characters; flipping the cell
# style back to Normal afterwards avoids leaving a stray style index on the
# cell (it keeps the unstyled look the original cells had).
function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '30.646.06'
Set-TextValue 'E2' '  +0.38%  '
Set-TextValue 'D3' '1.962.55'
Set-TextValue 'E3' '  +2.30%  '
Set-TextValue 'D4' '0.9997'
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '249.95'
Set-TextValue 'E5' '  +1.80%  '
Set-TextValue 'D6' '0.9998'
Set-TextValue 'E6' '  -0.02%  '
Set-TextValue 'D7' '0.4845'
Set-TextValue 'E7' '  +0.99%  '
Set-TextValue 'D8' '0.2944'
Set-TextValue 'E8' '  +1.60%  '
Set-TextValue 'D9' '0.06781'
Set-TextValue 'E9' '  +0.89%  '
Set-TextValue 'D10' '109.85'
Set-TextValue 'E10' '  -0.80%  '
Set-TextValue 'D11' '19.40'
Set-TextValue 'E11' '  +1.94%  '
Set-TextValue 'D12' '1.963.01'
Set-TextValue 'E12' '  +2.39%  '
Set-TextValue 'D13' '0.07762'
Set-TextValue 'E13' '  +2.56%  '
Set-TextValue 'D14' '5.457'
Set-TextValue 'E14' '  +3.51%  '
Set-TextValue 'D15' '0.6882'
Set-TextValue 'E15' '  +3.07%  '
Set-TextValue 'D16' '294.58'
Set-TextValue 'E16' '  -1.62%  '
Set-TextValue 'D17' '30.662.50'
Set-TextValue 'E17' '  +0.45%  '
Set-TextValue 'D18' '13.24'
Set-TextValue 'E18' '  +1.83%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D19' '0.000007705'
Set-TextValue 'E19' '  +1.70%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D20' '2.223.58'
Set-TextValue 'E20' '  +2.66%  '
Set-TextValue 'D21' '5.614'
Set-TextValue 'E21' '  -0.04%  '
Set-TextValue 'D22' '0.9998'
Set-TextValue 'E22' '  -0.03%  '
Set-TextValue 'D23' '0.9997'
Set-TextValue 'E23' '  -0.10%  '
Set-TextValue 'D24' '6.597'
Set-TextValue 'E24' '  +1.66%  '
Set-TextValue 'D25' '9.886'
Set-TextValue 'E25' '  +4.17%  '
Set-TextValue 'D26' '170.24'
Set-TextValue 'E26' '  +3.41%  '
Set-TextValue 'D27' '20.16'
Set-TextValue 'E27' '  -0.58%  '
Set-TextValue 'D28' '2.180'
Set-TextValue 'E28' '  +3.13%  '
Set-TextValue 'D30' '1.435'
Set-TextValue 'E30' '  +2.55%  '
Set-TextValue 'D31' '4.714'
Set-TextValue 'E31' '  +16.45%  '
Set-TextValue 'D32' '4.441'
Set-TextValue 'E32' '  +6.54%  '
Set-TextValue 'D33' '0.05112'
Set-TextValue 'E33' '  +2.21%  '
Set-TextValue 'D34' '0.7700'
Set-TextValue 'E34' '  +4.41%  '
Set-TextValue 'D35' '1.181'
Set-TextValue 'E35' '  +3.92%  '
Set-TextValue 'E36' '  +0.64%  '
Set-TextValue 'D37' '2.733'
Set-TextValue 'E37' '  +0.42%  '
Set-TextValue 'D38' '2.721'
Set-TextValue 'E38' '  +1.34%  '
Set-TextValue 'D39' '2.126'
Set-TextValue 'E39' '  +5.12%  '
Set-TextValue 'D40' '6.422'
Set-TextValue 'E40' '  +9.24%  '
Set-TextValue 'D41' '0.4474'
Set-TextValue 'E41' '  +0.97%  '
Set-TextValue 'D42' '109.11'
Set-TextValue 'D43' '0.8765'
Set-TextValue 'E43' '  +1.63%  '
Set-TextValue 'D44' '70.28'
Set-TextValue 'E44' '  -2.71%  '
Set-TextValue 'E45' '  +0.06%  '
Set-TextValue 'D46' '7.474'
Set-TextValue 'E46' '  +2.67%  '
Set-TextValue 'D47' '0.1283'
Set-TextValue 'E47' '  +4.22%  '
Set-TextValue 'D48' '9.360'
Set-TextValue 'E48' '  +0.53%  '
Set-TextValue 'D49' '35.97'
Set-TextValue 'E49' '  +2.51%  '
Set-TextValue 'D50' '47.62'
Set-TextValue 'E50' '  -3.48%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D51' '0.4072'
Set-TextValue 'E51' '  +1.50%  '
